# Update odds values for rows 2, 5, and 10 on Sheet1
# to reflect refreshed FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 3
$ws.Range("N2").Value = 8
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 1.83
$ws.Range("X2").Value = 10
$ws.Range("AC2").Value = 8
$ws.Range("AG2").Value = 301
$ws.Range("AK2").Value = 34
$ws.Range("AN2").Value = 4.33
$ws.Range("AO2").Value = 13
$ws.Range("AT2").Value = 2.63
$ws.Range("AV2").Value = 51
$ws.Range("AX2").Value = 19
$ws.Range("BA2").Value = 81
$ws.Range("G5").Value = 3.6
$ws.Range("H5").Value = 2.9
$ws.Range("I5").Value = 2.25
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 1.95
$ws.Range("L5").Value = 3.1
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("W5").Value = 8.5
$ws.Range("AD5").Value = 5.5
$ws.Range("AI5").Value = 10
$ws.Range("AK5").Value = 21
$ws.Range("AL5").Value = 21
$ws.Range("AN5").Value = 5
$ws.Range("AQ5").Value = 67
$ws.Range("AT5").Value = 2.38
$ws.Range("AX5").Value = 13
$ws.Range("BA5").Value = 81
$ws.Range("BB5").Value = 251
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
